$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 585, shifting existing rows 585-626 down to 586-627
$ws.Rows.Item(585).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A585").Value = 8
$ws.Range("B585").Value = "Terminal La Palmera de La Serena"
$ws.Range("C585").Value = "Coquimbo"
$ws.Range("D585").Value = 44746
$ws.Range("E585").Value = 4
$ws.Range("F585").Value = 100112024
$ws.Range("G585").Value = "Choclo"
$ws.Range("H585").Value = "Dulce o Americano"
$ws.Range("I585").Value = "Primera"
$ws.Range("J585").Value = 480
$ws.Range("K585").Value = 43000
$ws.Range("L585").Value = 44000
$ws.Range("M585").Value = 43500
$ws.Range("N585").Value = "$/malla 70 unidades"
$ws.Range("O585").Value = "Región de Arica y Parinacota"
$ws.Range("P585").Value = 621
$ws.Range("Q585").Value = 70
$ws.Range("R585").Value = "Hortaliza"
